$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.738.30"
$ws.Range("E2").Value = "  +2.95%  "
$ws.Range("D3").Value = "2.603.55"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.600"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("E9").Value = "  +1.91%  "
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("E11").Value = "  +2.77%  "
$ws.Range("E12").Value = "  -4.48%  "
$ws.Range("E13").Value = "  +7.09%  "
$ws.Range("D14").Value = "3.065.27"
$ws.Range("E14").Value = "  +1.47%  "
$ws.Range("D15").Value = "60.751.80"
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.36%  "
$ws.Range("E17").Value = "  +3.01%  "
$ws.Range("D18").Value = "2.613.53"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +13.09%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("E24").Value = "  +13.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.13"
$ws.Range("D25").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.65%  "
$ws.Range("D29").Value = "0.0₃0791"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("E30").Value = "  +9.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.85%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "162.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.24"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.959"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.85%  "
$ws.Range("E37").Value = "  +3.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +7.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "37.72"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("E40").Value = "  -2.58%  "
$ws.Range("E41").Value = "  +3.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "296.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.99"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("E47").Value = "  +1.84%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0241"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.65%  "
